$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3577081263065338
$ws.Range("B1").Value = 0.2669700682163239
$ws.Range("C1").Value = 1.097843289375305
$ws.Range("D1").Value = 3.678993225097656
$ws.Range("E1").Value = 1.944491386413574
